# Add Flow vs R1L summary block (rows 38-40, columns G:J) to Sheet1,
# mirroring the existing Kpl summary block at rows 22-24 (columns C:F)
# but built from column F ("Flow_Lac") instead of column B ("Kpl").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 38): same column headers as row 22 (HK-2, UMRC6, UOK262, UOK + DIDS)
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# Row label (row 39): new series name
$ws.Range("F39").Value = "Flow_Lac"

# Row 39: averages of column F over the same row groupings used for column B
$ws.Range("G39").Formula = "=AVERAGE(F`$1:F`$3)"
$ws.Range("H39").Formula = "=AVERAGE(F`$4:F`$6)"
$ws.Range("I39").Formula = "=AVERAGE(F`$9:F`$11)"
$ws.Range("J39").Formula = "=AVERAGE(F`$13:F`$16)"

# Row 40: standard error of the mean (STDEV/SQRT(COUNT)) over the same groupings
$ws.Range("G40").Formula = "=STDEV(F`$1:F`$3)/SQRT(COUNT(F`$1:F`$3))"
$ws.Range("H40").Formula = "=STDEV(F`$4:F`$6)/SQRT(COUNT(F`$4:F`$6))"
$ws.Range("I40").Formula = "=STDEV(F`$9:F`$11)/SQRT(COUNT(F`$9:F`$11))"
$ws.Range("J40").Formula = "=STDEV(F`$13:F`$16)/SQRT(COUNT(F`$13:F`$16))"

# Update selection to match the saved view (cell F10)
$ws.Range("F10").Select()
